$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 54 data
# Copy the style of A53 (date-formatted cell) onto A54 so the same style index is reused
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(54, 1).Value = 43569
$ws.Cells.Item(54, 2).Value = 34
$ws.Cells.Item(54, 3).Value = 477
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 5).Value = 0
$ws.Cells.Item(54, 6).Value = 463
$ws.Cells.Item(54, 7).Value = 225
$ws.Cells.Item(54, 8).Value = 252

# Update selection to match diff: the whole new row 54 ends up selected
# (activeCell moves onto row 54; sqref covers the entire row A54:XFD54)
$ws.Rows("54:54").Select()

# Keep the view scrolled to where it was (topLeftCell A16)
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
